$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 39888
$ws.Range("J108").Value = 39888
$ws.Range("L108").Value = 39888
$ws.Range("N108").Value = -47568
$ws.Range("H116").Value = 4899.35
$ws.Range("I116").Value = 2363.5454
$ws.Range("J116").Value = 7998.6665
$ws.Range("K116").Value = 2363.5454
$ws.Range("L116").Value = 7998.6665
$ws.Range("M116").Value = 1078.4546
$ws.Range("N116").Value = -14882.6665
$ws.Range("H123").Value = 39297.586
$ws.Range("J123").Value = 39297.586
$ws.Range("L123").Value = 39297.586
$ws.Range("N123").Value = -49097.586
$ws.Range("H124").Value = 47607.5
$ws.Range("J124").Value = 47607.5
$ws.Range("L124").Value = 47607.5
$ws.Range("N124").Value = -57427.5
$ws.Range("H126").Value = 54091.668
$ws.Range("J126").Value = 54091.668
$ws.Range("L126").Value = 54091.668
$ws.Range("N126").Value = -63971.668
$ws.Range("H128").Value = 34693
$ws.Range("J128").Value = 34693
$ws.Range("L128").Value = 34693
$ws.Range("N128").Value = -44653

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20639.096
$ws.Range("I32").Value = 4295.091
$ws.Range("J32").Value = 133004.12
$ws.Range("K32").Value = 4295.091
$ws.Range("L32").Value = 133004.12
$ws.Range("M32").Value = -4008.091
$ws.Range("N32").Value = -133578.12
$ws.Range("H103").Value = 39478
$ws.Range("J103").Value = 39478
$ws.Range("L103").Value = 39478
$ws.Range("N103").Value = -41822
$ws.Range("H122").Value = 2300.8
$ws.Range("I122").Value = 1964.7273
$ws.Range("J122").Value = 3225
$ws.Range("K122").Value = 5894.1819
$ws.Range("L122").Value = 9675
$ws.Range("M122").Value = -3444.1819
$ws.Range("N122").Value = -14575
$ws.Range("H123").Value = 1275000
$ws.Range("J123").Value = 1275000
$ws.Range("L123").Value = 1275000
$ws.Range("N123").Value = -1284800
$ws.Range("H125").Value = 150049920
$ws.Range("J125").Value = 150049920
$ws.Range("L125").Value = 150049920
$ws.Range("N125").Value = -150059760
$ws.Range("H128").Value = 56392
$ws.Range("J128").Value = 56392
$ws.Range("L128").Value = 56392
$ws.Range("N128").Value = -66352
$ws.Range("H129").Value = 43331.332
$ws.Range("J129").Value = 43997.6
$ws.Range("L129").Value = 43997.6
$ws.Range("N129").Value = -53997.6
$ws.Range("H133").Value = 52236.2
$ws.Range("J133").Value = 52236.2
$ws.Range("L133").Value = 52236.2
$ws.Range("N133").Value = -57296.2

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 37680
$ws.Range("J124").Value = 37680
$ws.Range("L124").Value = 37680
$ws.Range("N124").Value = -47500
$ws.Range("H125").Value = 53090
$ws.Range("J125").Value = 53090
$ws.Range("L125").Value = 53090
$ws.Range("N125").Value = -62930
$ws.Range("H132").Value = 55000
$ws.Range("J132").Value = 55000
$ws.Range("L132").Value = 55000
$ws.Range("N132").Value = -65120

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 17214.5
$ws.Range("J97").Value = 17214.5
$ws.Range("L97").Value = 17214.5
$ws.Range("N97").Value = -19196.5
$ws.Range("H123").Value = 39818.57
$ws.Range("J123").Value = 39818.57
$ws.Range("L123").Value = 39818.57
$ws.Range("N123").Value = -49618.57
$ws.Range("H130").Value = 35976
$ws.Range("J130").Value = 35976
$ws.Range("L130").Value = 35976
$ws.Range("N130").Value = -46016

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 19077.5
$ws.Range("J62").Value = 19077.5
$ws.Range("L62").Value = 19077.5
$ws.Range("N62").Value = -20449.5
$ws.Range("H65").Value = 19077.5
$ws.Range("J65").Value = 19077.5
$ws.Range("L65").Value = 57232.5
$ws.Range("N65").Value = -64096.5
$ws.Range("H128").Value = 51714
$ws.Range("J128").Value = 51714
$ws.Range("L128").Value = 51714
$ws.Range("N128").Value = -61674
$ws.Range("H130").Value = 55333.75
$ws.Range("J130").Value = 55333.75
$ws.Range("L130").Value = 55333.75
$ws.Range("N130").Value = -65373.75
$ws.Range("H133").Value = 38500
$ws.Range("J133").Value = 38500
$ws.Range("L133").Value = 38500
$ws.Range("N133").Value = -48620
$ws.Range("H135").Value = 55457.5
$ws.Range("J135").Value = 55457.5
$ws.Range("L135").Value = 55457.5
$ws.Range("N135").Value = -65597.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 35500
$ws.Range("J80").Value = 35500
$ws.Range("L80").Value = 35500
$ws.Range("N80").Value = -37746
$ws.Range("H83").Value = 35500
$ws.Range("J83").Value = 35500
$ws.Range("L83").Value = 106500
$ws.Range("N83").Value = -117732
$ws.Range("H108").Value = 24163.334
$ws.Range("J108").Value = 24163.334
$ws.Range("L108").Value = 24163.334
$ws.Range("N108").Value = -31843.334
$ws.Range("H125").Value = 44450
$ws.Range("J125").Value = 44450
$ws.Range("L125").Value = 44450
$ws.Range("N125").Value = -54290
$ws.Range("H127").Value = 50418.57
$ws.Range("J127").Value = 50418.57
$ws.Range("L127").Value = 50418.57
$ws.Range("N127").Value = -60338.57
$ws.Range("H128").Value = 50182.375
$ws.Range("J128").Value = 50182.375
$ws.Range("L128").Value = 50182.375
$ws.Range("N128").Value = -60142.375
$ws.Range("H129").Value = 32847.6
$ws.Range("J129").Value = 32847.6
$ws.Range("L129").Value = 32847.6
$ws.Range("N129").Value = -42847.6
$ws.Range("H132").Value = 2391.4194
$ws.Range("I132").Value = 1690.7333
$ws.Range("J132").Value = 3048.3125
$ws.Range("K132").Value = 5072.199900000001
$ws.Range("L132").Value = 9144.9375
$ws.Range("M132").Value = -2542.199900000001
$ws.Range("N132").Value = -14204.9375
$ws.Range("H133").Value = 52631.75
$ws.Range("J133").Value = 52631.75
$ws.Range("L133").Value = 52631.75
$ws.Range("N133").Value = -57691.75
$ws.Range("H134").Value = 55714.5
$ws.Range("J134").Value = 55714.5
$ws.Range("L134").Value = 55714.5
$ws.Range("N134").Value = -65854.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 17195.4
$ws.Range("J21").Value = 17195.4
$ws.Range("L21").Value = 17195.4
$ws.Range("N21").Value = -17665.4
$ws.Range("H35").Value = 17195.4
$ws.Range("J35").Value = 17195.4
$ws.Range("L35").Value = 17195.4
$ws.Range("N35").Value = -17775.4
$ws.Range("H54").Value = 7261.6
$ws.Range("J54").Value = 7261.6
$ws.Range("L54").Value = 7261.6
$ws.Range("N54").Value = -8301.6
$ws.Range("H64").Value = 24620.445
$ws.Range("J64").Value = 24620.445
$ws.Range("L64").Value = 24620.445
$ws.Range("N64").Value = -25116.445
$ws.Range("H67").Value = 24620.445
$ws.Range("J67").Value = 24620.445
$ws.Range("L67").Value = 24620.445
$ws.Range("N67").Value = -26336.445
$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33622
$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -104112
$ws.Range("H125").Value = 47967.855
$ws.Range("J125").Value = 47967.855
$ws.Range("L125").Value = 47967.855
$ws.Range("N125").Value = -57807.855
$ws.Range("H127").Value = 51108.375
$ws.Range("J127").Value = 51108.375
$ws.Range("L127").Value = 51108.375
$ws.Range("N127").Value = -61028.375
$ws.Range("H129").Value = 25762.334
$ws.Range("J129").Value = 25762.334
$ws.Range("L129").Value = 25762.334
$ws.Range("N129").Value = -35762.334
$ws.Range("H130").Value = 28627.25
$ws.Range("J130").Value = 28627.25
$ws.Range("L130").Value = 28627.25
$ws.Range("N130").Value = -38667.25
